# resolver errores de tipo_equipo y resolver problema de tablas en sql con
# mayusculas no funcionando en linux
#
# Sheet "AIO": rewrite rows 2 & 3 with corrected values, drop the now
# duplicate/bogus rows 4-6 entirely.
# Sheet "Otros": drop rows 2-5, keep only the header row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "AIO"
# ---------------------------------------------------------------------
$aio = $wb.Worksheets.Item("AIO")

# Row 2
# (D2/J2/K2 hold the digit-string "123" - prefix with an apostrophe so it
# is stored as text instead of being auto-converted to a number.)
$aio.Range("A2").Value = "Concepcion"
$aio.Range("B2").Value = "Coronel"
$aio.Range("C2").Value = "OFICINA"
$aio.Range("D2").Value = "'123"
$aio.Range("E2").Value = "BALDOMERO LILLO"
$aio.Range("F2").Value = 8015892
$aio.Range("G2").Value = "AIO"
$aio.Range("H2").Value = "Lenovo"
$aio.Range("I2").Value = "V2414"
$aio.Range("J2").Value = "'123"
$aio.Range("K2").Value = "'123"
$aio.Range("L2").Value = "Sonda"

# Row 3
$aio.Range("A3").Value = "Concepcion"
$aio.Range("B3").Value = "Coronel"
$aio.Range("C3").Value = "OFICINA"
$aio.Range("D3").Value = "'123"
$aio.Range("E3").Value = "BALDOMERO LILLO"
$aio.Range("F3").Value = 8015892
$aio.Range("G3").Value = "AIO"
$aio.Range("H3").Value = "Lenovo"
$aio.Range("I3").Value = "V2414"
$aio.Range("J3").Value = "'123"
$aio.Range("K3").Value = "'123"
$aio.Range("L3").Value = "Sonda"

# Rows 4-6 are no longer needed, drop them (shifts nothing else up since
# they are the last rows).
$aio.Range("A4:A6").EntireRow.Delete()

# ---------------------------------------------------------------------
# Sheet "Otros"
# ---------------------------------------------------------------------
$otros = $wb.Worksheets.Item("Otros")

# Rows 2-5 go away, only the header row is kept.
$otros.Range("A2:A5").EntireRow.Delete()
